$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 65, pushing the existing rows 65-96 down to 66-97
$ws.Rows.Item(65).Insert()

# Populate the newly inserted row 65 with the new weekly price record
$ws.Range("A65").Value = 5
$ws.Range("B65").Value = "Macroferia Regional de Talca"
$ws.Range("C65").Value = "Maule"
$ws.Range("D65").Value = 44488
$ws.Range("E65").Value = 7
$ws.Range("F65").Value = 100112031
$ws.Range("G65").Value = "Poroto verde"
$ws.Range("H65").Value = "Sin especificar"
$ws.Range("I65").Value = "Primera"
$ws.Range("J65").Value = 100
$ws.Range("K65").Value = 40000
$ws.Range("L65").Value = 40000
$ws.Range("M65").Value = 40000
$ws.Range("N65").Value = "`$/malla 25 kilos"
$ws.Range("O65").Value = "Región de Arica y Parinacota"
$ws.Range("P65").Value = 1600
$ws.Range("Q65").Value = 25
$ws.Range("R65").Value = "Hortaliza"
